$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top; existing rows shift down by one.
$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = "from"
$ws.Range("B1").Value = "to"
$ws.Range("C1").Value = "cap"
$ws.Range("D1").Value = "dis"
